$d = $word.ActiveDocument

$replacements = @(
    @("592÷5=118, 2", "399÷5=79, 4"),
    @("816÷2=408, 0", "654÷6=109, 0"),
    @("231÷8=28, 7", "215÷4=53, 3"),
    @("132÷3=44, 0", "219÷2=109, 1"),
    @("145÷8=18, 1", "256÷5=51, 1"),
    @("963÷5=192, 3", "497÷8=62, 1"),
    @("878÷9=97, 5", "876÷9=97, 3"),
    @("738÷3=246, 0", "777÷6=129, 3"),
    @("435÷6=72, 3", "385÷4=96, 1"),
    @("484÷9=53, 7", "561÷2=280, 1"),
    @("786÷9=87, 3", "348÷4=87, 0"),
    @("473÷9=52, 5", "342÷2=171, 0"),
    @("408÷9=45, 3", "774÷3=258, 0"),
    @("599÷5=119, 4", "511÷9=56, 7"),
    @("558÷5=111, 3", "552÷8=69, 0"),
    @("555÷4=138, 3", "464÷2=232, 0"),
    @("227÷4=56, 3", "624÷8=78, 0"),
    @("817÷4=204, 1", "564÷4=141, 0"),
    @("164÷4=41, 0", "372÷7=53, 1"),
    @("608÷8=76, 0", "972÷9=108, 0"),
    @("567÷5=113, 2", "679÷7=97, 0"),
    @("459÷5=91, 4", "362÷4=90, 2"),
    @("676÷2=338, 0", "924÷5=184, 4"),
    @("730÷4=182, 2", "732÷2=366, 0"),
    @("823÷2=411, 1", "410÷4=102, 2")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
